$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 1.919867272924993

# Row 3
$ws.Range("B3").Value = 0.3048080303191223
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 27.33178407655716

# Row 4
$ws.Range("B4").Value = 0.04763786555579896
$ws.Range("C4").Value = 0.04240448674262143
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 1.392311199053152

# Row 5
$ws.Range("B5").Value = 3.230985683306322
$ws.Range("C5").Value = 114.8270160096505
$ws.Range("D5").Value = 0.1575252929769615
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("G5").Value = 126.8757594718828
